# The author removed every slide from the deck (the slide list collapses to
# an empty <p:sldIdLst>, which PowerPoint then omits entirely). Deleting a
# slide that is the sole user of its notes page also drops the now-orphaned
# notes slide, matching the upstream diff (ppt/notesSlides/notesSlide1.xml,
# ppt/slides/slide1.xml, slide2.xml and slide3.xml all disappear while the
# slide master / layouts / theme stay untouched).

$p = $ppt.ActivePresentation

for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $p.Slides.Item($i).Delete()
}
